# KNIPS Main input data.xlsx - align shared string / column names with readxl
# package expectations: disambiguate generic "ll"/"UL" CI-bound labels and the
# generic "xb:"/"dxb:" covariance row/column labels so that readxl (which
# requires unique column names) can parse the sheets correctly.

$wb = $excel.ActiveWorkbook

# --- utilities sheet: give the confidence-interval columns unique names ---
$wsUtil = $wb.Worksheets.Item("utilities")
$wsUtil.Range("E1").Value = "pre primary ll"
$wsUtil.Range("F1").Value = "pre primary UL"
$wsUtil.Range("I1").Value = "6 months after primary ll"
$wsUtil.Range("J1").Value = "6 months after primary UL"
$wsUtil.Range("M1").Value = "pre revision ll"
$wsUtil.Range("N1").Value = "pre revision UL"
$wsUtil.Range("Q1").Value = "6 months after revision ll"
$wsUtil.Range("R1").Value = "6 months after revision UL"
$wsUtil.Range("R2").Select()

# --- second_revision_covariance sheet: label the rows/columns explicitly ---
$wsCov = $wb.Worksheets.Item("second_revision_covariance")
$wsCov.Range("A1").Value = "names"
$wsCov.Range("B1").Value = "xb1"
$wsCov.Range("C1").Value = "xb2"
$wsCov.Range("D1").Value = "xb3"
$wsCov.Range("E1").Value = "xb4"
$wsCov.Range("F1").Value = "xb5"
$wsCov.Range("G1").Value = "xb6"
$wsCov.Range("H1").Value = "xb7"
$wsCov.Range("I1").Value = "xb8"
$wsCov.Range("B1:I1").Select()

# --- restore the originally active sheet/tab ---
$wsActive = $wb.Worksheets.Item("90d_mortality")
$wsActive.Activate()
